$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 into I1:J1, then set header labels
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate I0/IF values for rows 2-69
$data = @(
    @(2, 6, 6),
    @(3, 7, 8),
    @(4, 8, 9),
    @(5, 7, 7),
    @(6, 8, 9),
    @(7, 8, 8),
    @(8, 9, 9),
    @(9, 6, 7),
    @(10, 8, 8),
    @(11, 8, 8),
    @(12, 9, 9),
    @(13, 9, 9),
    @(14, 7, 7),
    @(15, 8, 8),
    @(16, 8, 8),
    @(17, 8, 8),
    @(18, 8, 8),
    @(19, 9, 9),
    @(20, 9, 9),
    @(21, 9, 9),
    @(22, 9, 9),
    @(23, 9, 9),
    @(24, 8, 8),
    @(25, 8, 8),
    @(26, 9, 9),
    @(27, 9, 9),
    @(28, 9, 9),
    @(29, 9, 9),
    @(30, 9, 9),
    @(31, 8, 8),
    @(32, 8, 8),
    @(33, 9, 9),
    @(34, 8, 9),
    @(35, 8, 8),
    @(36, 9, 9),
    @(37, 8, 8),
    @(38, 10, 10),
    @(39, 8, 8),
    @(40, 7, 7),
    @(41, 9, 9),
    @(42, 7, 7),
    @(43, 7, 7),
    @(44, 8, 8),
    @(45, 8, 9),
    @(46, 8, 8),
    @(47, 8, 8),
    @(48, 8, 8),
    @(49, 8, 8),
    @(50, 6, 7),
    @(51, 8, 8),
    @(52, 7, 7),
    @(53, 5, 6),
    @(54, 7, 7),
    @(55, 6, 6),
    @(56, 7, 7),
    @(57, 7, 7),
    @(58, 6, 7),
    @(59, 5, 6),
    @(60, 7, 7),
    @(61, 6, 6),
    @(62, 4, 5),
    @(63, 7, 8),
    @(64, 8, 8),
    @(65, 7, 7),
    @(66, 7, 7),
    @(67, 8, 8),
    @(68, 7, 7),
    @(69, 4, 4),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
